$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Mã giáo vụ"
$ws.Range("B1").Value = "Họ và tên"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Viện"
$ws.Range("E1").Value = "Khóa công khai"

# Row 2 - Nguyễn Văn B
$ws.Range("A2").Value = "GV1234"
$ws.Range("B2").Value = "Nguyễn Văn B"
$ws.Range("C2").Value = "nguyenvanb@soict.hust.edu.vn"
$ws.Range("D2").Value = "Viện CNTT&TT"
$ws.Range("E2").Value = "31d0e835390695f825a1322b38bdb3de71c075..."

# Row 3 - Lý Thị C
$ws.Range("A3").Value = "GV1235"
$ws.Range("B3").Value = "Lý Thị C"
$ws.Range("C3").Value = "lythic@spkt.hust.edu.vn"
$ws.Range("D3").Value = "Viện Sư phạm Kỹ thuật"
$ws.Range("E3").Value = "81741c77ad5e5ff27ec91a94ced51b82a37968..."

# Row 4 - Lê Thị D
$ws.Range("A4").Value = "GV1236"
$ws.Range("B4").Value = "Lê Thị D"
$ws.Range("C4").Value = "lethidc@nn.hust.edu.vn"
$ws.Range("D4").Value = "Viện Ngoại ngữ"
$ws.Range("E4").Value = "ecfdc8f69b08d0260ba2309d7b8e064a28f0eb5..."

# Row 5 - Trần Văn E
$ws.Range("A5").Value = "GV1237"
$ws.Range("B5").Value = "Trần Văn E"
$ws.Range("C5").Value = "tranvane@dtvt.hust.edu.vn"
$ws.Range("D5").Value = "Viện Điện tử viễn thông"
$ws.Range("E5").Value = "4d4ebfbf5ea1f3b61b04434528844956ab6890536..."

# Row 6 - Đào Thị F
$ws.Range("A6").Value = "GV1238"
$ws.Range("B6").Value = "Đào Thị F"
$ws.Range("C6").Value = "daothif@dktdh.hust.edu.vn"
$ws.Range("D6").Value = "Viện Điều khiển Tự động hóa"
$ws.Range("E6").Value = "4d4ebfbf5ea1f3b61b04434528844956ab6890..."

# Re-fit the data columns (B:E) to their new content widths
$ws.Columns("B:B").ColumnWidth = 11.25
$ws.Columns("C:C").ColumnWidth = 25.9375
$ws.Columns("D:D").ColumnWidth = 23.59375
$ws.Columns("E:E").ColumnWidth = 42.65625

# Selection moved
[void]$ws.Range("C15").Select()

# Page setup: portrait orientation
$ws.PageSetup.Orientation = 1
